$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns stay text-formatted, matching the
# original inlineStr cells (avoids "1.00" -> 1 numeric coercion).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "67.314.82"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "2.624.52"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "597.46"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "168.19"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("D9").Value = "2.623.19"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("D13").Value = "5.24"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "27.79"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").Value = "3.103.85"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "67.347.99"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "2.624.61"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("E19").Value = "  +4.46%  "
$ws.Range("D20").Value = "8.06"
$ws.Range("E20").Value = "  +7.53%  "
$ws.Range("D21").Value = "356.73"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("D24").Value = "10.45"
$ws.Range("E24").Value = "  +3.84%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  -4.79%  "
$ws.Range("D27").Value = "69.66"
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("D28").Value = "2.761.07"
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("D31").Value = "550.67"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").Value = "7.94"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").Value = "0.135"
$ws.Range("E35").Value = "  +5.25%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("D37").Value = "1.51"
$ws.Range("E37").Value = "  -3.74%  "
$ws.Range("D38").Value = "158.38"
$ws.Range("E38").Value = "  +2.72%  "
$ws.Range("D39").Value = "19.00"
$ws.Range("E39").Value = "  -2.34%  "
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("E45").Value = "  -3.83%  "
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "0.0₆0295"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").Value = "152.05"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("D50").Value = "3.78"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("E51").Value = "  -1.20%  "
